$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.513.23"
$ws.Range("E2").Value = "'  -1.91%  "

$ws.Range("D3").Value = "'1.964.32"
$ws.Range("E3").Value = "'  +0.11%  "

$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "'  +0.01%  "

$ws.Range("D5").Value = "'322.68"
$ws.Range("E5").Value = "'  -1.45%  "

$ws.Range("D6").Value = "'1.011"
$ws.Range("E6").Value = "'  -0.04%  "

$ws.Range("D7").Value = "'0.4763"
$ws.Range("E7").Value = "'  -4.20%  "

$ws.Range("D8").Value = "'0.4051"
$ws.Range("E8").Value = "'  -3.67%  "

$ws.Range("D9").Value = "'54.04"
$ws.Range("E9").Value = "'  -0.19%  "

$ws.Range("D10").Value = "'0.08480"
$ws.Range("E10").Value = "'  -5.99%  "

$ws.Range("E11").Value = "'  -3.27%  "

$ws.Range("D12").Value = "'22.42"
$ws.Range("E12").Value = "'  -2.27%  "

$ws.Range("D13").Value = "'1.959.47"
$ws.Range("E13").Value = "'  +0.04%  "

$ws.Range("D14").Value = "'7.604"
$ws.Range("E14").Value = "'  -3.03%  "

$ws.Range("D15").Value = "'6.189"
$ws.Range("E15").Value = "'  -3.52%  "

$ws.Range("D16").Value = "'1.013"
$ws.Range("E16").Value = "'  +0.07%  "

$ws.Range("D17").Value = "'90.96"
$ws.Range("E17").Value = "'  +0.01%  "

$ws.Range("E18").Value = "'  -2.34%  "

$ws.Range("D19").Value = "'0.06642"
$ws.Range("E19").Value = "'  -0.17%  "

$ws.Range("D20").Value = "'18.53"
$ws.Range("E20").Value = "'  -3.27%  "

$ws.Range("E21").Value = "'  +0.13%  "

$ws.Range("D22").Value = "'5.866"
$ws.Range("E22").Value = "'  -1.08%  "

$ws.Range("D23").Value = "'28.561.08"
$ws.Range("E23").Value = "'  -1.82%  "

$ws.Range("E24").Value = "'  -3.72%  "

$ws.Range("D25").Value = "'2.300"
$ws.Range("E25").Value = "'  +0.35%  "

$ws.Range("D26").Value = "'2.235.41"
$ws.Range("E26").Value = "'  +1.93%  "

$ws.Range("D27").Value = "'155.85"
$ws.Range("E27").Value = "'  -0.31%  "

$ws.Range("D28").Value = "'20.36"
$ws.Range("E28").Value = "'  -1.09%  "

$ws.Range("E29").Value = "'  -5.21%  "

$ws.Range("D30").Value = "'2.168"
$ws.Range("E30").Value = "'  -3.80%  "

$ws.Range("D31").Value = "'124.71"
$ws.Range("E31").Value = "'  -1.65%  "

$ws.Range("D32").Value = "'0.9826"
$ws.Range("E32").Value = "'  -5.29%  "

$ws.Range("D33").Value = "'0.09639"
$ws.Range("E33").Value = "'  -1.97%  "

$ws.Range("D34").Value = "'1.457"
$ws.Range("E34").Value = "'  -4.21%  "

$ws.Range("D35").Value = "'3.701"
$ws.Range("E35").Value = "'  -0.02%  "

$ws.Range("D36").Value = "'5.633"
$ws.Range("E36").Value = "'  -2.93%  "

$ws.Range("D37").Value = "'9.120"
$ws.Range("E37").Value = "'  +1.40%  "

$ws.Range("D38").Value = "'0.02332"
$ws.Range("E38").Value = "'  -3.67%  "

$ws.Range("D39").Value = "'0.06231"
$ws.Range("E39").Value = "'  -1.45%  "

$ws.Range("D40").Value = "'1.254"
$ws.Range("E40").Value = "'  -2.45%  "

$ws.Range("D41").Value = "'0.6213"
$ws.Range("E41").Value = "'  -3.24%  "

$ws.Range("D42").Value = "'11.18"
$ws.Range("E42").Value = "'  -2.29%  "

$ws.Range("D43").Value = "'1.011"
$ws.Range("E43").Value = "'  +0.07%  "

$ws.Range("D44").Value = "'0.1919"
$ws.Range("E44").Value = "'  -3.52%  "

$ws.Range("D45").Value = "'1.347"
$ws.Range("E45").Value = "'  +4.85%  "

$ws.Range("B46").Value = "'Decentraland"
$ws.Range("C46").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5952"
$ws.Range("E46").Value = "'  -3.73%  "

$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'13.02"
$ws.Range("E47").Value = "'  -3.27%  "

$ws.Range("D48").Value = "'2.061"
$ws.Range("E48").Value = "'  -4.95%  "

$ws.Range("D49").Value = "'3.410"
$ws.Range("E49").Value = "'  -1.91%  "

$ws.Range("D50").Value = "'0.06815"
$ws.Range("E50").Value = "'  -0.74%  "

$ws.Range("D51").Value = "'111.49"
$ws.Range("E51").Value = "'  -0.96%  "
